# Add the "RES" (boundary) attack block for SEED 1337: eight new epsilon
# columns (AK:AR) mirroring the existing FGSM block (AC:AJ), headed by a
# merged "BOUNDARY" label in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header label -------------------------------------------------
# Merge the new header block first (same as the other attack headers), then
# clone the formatting of the neighbouring "FGSM" header block (AC1:AJ1)
# onto the new AK1:AR1 block so every cell shares its plain style, and
# finally overwrite the label text.
$ws.Range("AK1:AR1").Merge()
$ws.Range("AC1:AJ1").Copy()
$ws.Range("AK1:AR1").PasteSpecial(-4122)
$ws.Range("AK1").Value = "BOUNDARY"

# --- Row 2: epsilon values ------------------------------------------------
# Same epsilon list as every other attack block: 0.01 .. 0.20
$ws.Range("AC2:AJ2").Copy()
$ws.Range("AK2:AR2").PasteSpecial(-4122)
$ws.Range("AC2:AJ2").Copy()
$ws.Range("AK2:AR2").PasteSpecial(-4163)

# --- Data rows -------------------------------------------------------------
# Row 4: LSTM / MAE
$ws.Range("AK4").Value = 385.8309564018249
$ws.Range("AL4").Value = 391.1910715484619
$ws.Range("AM4").Value = 394.7090356699626
$ws.Range("AN4").Value = 397.0392293294271
$ws.Range("AO4").Value = 406.3723960558573
$ws.Range("AP4").Value = 421.9438964970906
$ws.Range("AQ4").Value = 463.9439326985677
$ws.Range("AR4").Value = 618.506653251648

# Row 5: LSTM / RMSE
$ws.Range("AK5").Value = 491.1585717792276
$ws.Range("AL5").Value = 495.0655865475306
$ws.Range("AM5").Value = 499.3598978283595
$ws.Range("AN5").Value = 496.5037183535468
$ws.Range("AO5").Value = 515.1999076481044
$ws.Range("AP5").Value = 534.2222562855211
$ws.Range("AQ5").Value = 614.9062889104446
$ws.Range("AR5").Value = 801.1559861854735

# Row 6: LSTM / SIM
$ws.Range("AK6").Value = 0.9991558283219409
$ws.Range("AL6").Value = 0.9991399461906083
$ws.Range("AM6").Value = 0.9991279610152928
$ws.Range("AN6").Value = 0.9991520262174486
$ws.Range("AO6").Value = 0.9990737667233376
$ws.Range("AP6").Value = 0.998981776193121
$ws.Range("AQ6").Value = 0.998659902738153
$ws.Range("AR6").Value = 0.997716921676422

# Row 7: RNN / MAE
$ws.Range("AK7").Value = 415.4717549069723
$ws.Range("AL7").Value = 418.0133662605286
$ws.Range("AM7").Value = 418.7792081324259
$ws.Range("AN7").Value = 434.4067901674907
$ws.Range("AO7").Value = 433.9814000384013
$ws.Range("AP7").Value = 445.1295617039999
$ws.Range("AQ7").Value = 513.1038867441813
$ws.Range("AR7").Value = 680.8991062355042

# Row 8: RNN / RMSE
$ws.Range("AK8").Value = 511.9682169789609
$ws.Range("AL8").Value = 518.2065764221599
$ws.Range("AM8").Value = 517.6951448638611
$ws.Range("AN8").Value = 535.4154830575009
$ws.Range("AO8").Value = 530.5251851732917
$ws.Range("AP8").Value = 542.1869801055376
$ws.Range("AQ8").Value = 638.0499407098257
$ws.Range("AR8").Value = 827.2670269163315

# Row 9: RNN / SIM
$ws.Range("AK9").Value = 0.9990133123323123
$ws.Range("AL9").Value = 0.9989780500616255
$ws.Range("AM9").Value = 0.9989949035212186
$ws.Range("AN9").Value = 0.9989230801493203
$ws.Range("AO9").Value = 0.9989373682869108
$ws.Range("AP9").Value = 0.9989059283125489
$ws.Range("AQ9").Value = 0.998442591664765
$ws.Range("AR9").Value = 0.9973749075293667

# Row 10: GRU / MAE
$ws.Range("AK10").Value = 315.4416551526388
$ws.Range("AL10").Value = 317.3214560890198
$ws.Range("AM10").Value = 337.0510022735596
$ws.Range("AN10").Value = 348.4658999951681
$ws.Range("AO10").Value = 350.0478318850199
$ws.Range("AP10").Value = 380.9835699017843
$ws.Range("AQ10").Value = 441.200032749176
$ws.Range("AR10").Value = 673.8200302346547

# Row 11: GRU / RMSE
$ws.Range("AK11").Value = 441.2778873238339
$ws.Range("AL11").Value = 440.7485394175382
$ws.Range("AM11").Value = 459.6649198111583
$ws.Range("AN11").Value = 470.0443937153195
$ws.Range("AO11").Value = 479.7561652518424
$ws.Range("AP11").Value = 503.700016540021
$ws.Range("AQ11").Value = 564.7376069564607
$ws.Range("AR11").Value = 844.3186495523017

# Row 12: GRU / SIM
$ws.Range("AK12").Value = 0.9992702074576406
$ws.Range("AL12").Value = 0.999273539112866
$ws.Range("AM12").Value = 0.9992016846953634
$ws.Range("AN12").Value = 0.9991556022067656
$ws.Range("AO12").Value = 0.9991253388025124
$ws.Range("AP12").Value = 0.9990233690505527
$ws.Range("AQ12").Value = 0.9987860596935169
$ws.Range("AR12").Value = 0.9971681726709546
